$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diary")

# --- Header block: Name + Team ---
$ws.Range("B3").Value = "Babar Ayan"
$ws.Range("D3").Value = "Team 1"

# --- Row 10: fix existing first diary entry ---
$ws.Range("B10").Value = 45013
$ws.Range("C10").Value = "Erste Infoveranstaltung - Einführung in das Projekt und Vorstellung der Teams sowie erste grobe Einteilung."

# --- Row 11: KW 12 ---
$ws.Range("A11").Value = "KW 12"
$ws.Range("B11").NumberFormat = "m/d/yyyy"
$ws.Range("B11").Value = 45020
$ws.Range("C11").Value = "Einarbeitung in Thema: - Integrated Gradients + weitere Algorithmen + Papers lesen. Aufsetzen der Captum GUI und Ausführen aller Algorithmen. Ausführen einiger Algorithmen in PyCharm. Besprechung der Inhalte und Aufgabenaufteilung mit Teamkollegen. "
$ws.Range("D11").Value = "Clustering Algorithmen Matrix und Repräsentationsschemas. Prinzipien und Ideen hinter verschiedenen Algorithmen. Wichtigkeit und Motivation von Explainable AI. "
$ws.Range("E11").Value = "Aktuelle Bücher und Videos zum Thema vorgeschlagen. Richtung der Projektarbeit diskutiert. Aktuellen Stand zu den erarbeiteten Themen geteilt. Idee: GRAD-CAM als neuen Algorithmus später in GUI hinzuzufügen."
$ws.Range("F11").Value = 14

# --- Row 12: KW 13 ---
$ws.Range("A12").Value = "KW 13"
$ws.Range("B12").NumberFormat = "m/d/yyyy"
$ws.Range("B12").Value = "11.04.2023 (Ostern)"
$ws.Range("C12").Value = "Einarbeitung in die Funktionsweise von verschiedenen Algorithmen auf Captum. TCAV Paper gelesen und Videos zur funktionsweise angeschaut. Überblick zu den verschiedenen Algorithmen in Captum  bekommen. Austausch mit der Gruppe zu den Möglichkeiten Captum einzusetzen."
$ws.Range("D12").Value = "Nutzung von Stride anstatt Pooling. Bedienung von Captum GUI"
$ws.Range("E12").Value = "Austausch mit der Gruppe über die Funktionsweise von Captum."
$ws.Range("F12").Value = 10

# --- Row 13: KW 14 ---
$ws.Range("A13").Value = "KW 14"
$ws.Range("B13").NumberFormat = "m/d/yyyy"
$ws.Range("B13").Value = 45034
$ws.Range("C13").Value = "TCAV wiederholt und versucht die Schritte nachzuvollziehen."
$ws.Range("D13").Value = "Idee hinter mathematische Erklärung und Validierung von TCAV verstanden. "
$ws.Range("E13").Value = "Über funktionsweise von TCAV mit der Gruppe ausgetauscht"
$ws.Range("F13").Value = 12

# --- Row 14: KW 15 ---
$ws.Range("A14").Value = "KW 15"
$ws.Range("B14").NumberFormat = "m/d/yyyy"
$ws.Range("B14").Value = 45041
$ws.Range("C14").Value = "Präsentationsfüllen mit groben Stichpunkten befüllt und Reihenfolge der Inhalten vorgegeben. "
$ws.Range("D14").Value = "Fortgeschrittene Architekturen in Deep Learning, speziell ResNets und Positional Encoding bei Transformers + Nutzung von Transformers anstatt CNNs für Bildklassifikation."
$ws.Range("E14").Value = "Gliederung der Präsentation"
$ws.Range("F14").Value = 8

# --- Row 15: KW 16 ---
$ws.Range("A15").Value = "KW 16"
$ws.Range("B15").NumberFormat = "m/d/yyyy"
$ws.Range("B15").Value = 45048
$ws.Range("C15").Value = "Paper Präsentation vorbereitet. Ausarbeitung der Stichpunkte in den Präsentationsfolien. Hinzufügen von passenden Bildern (+ Quellenangaben). Für meinen Teil der Paperpräsentation habe ich die Folien deutlich reduziert von 8 auf 3, um die 5 Minuten Vortragszeit einzuhalten. Ich habe speziell die Folieninhalte zu Einleitung, Saliency Maps, TCAV und Ziele gestaltet. Feedback und Unterstützung bei der Formulierung der Folien zur Validierung von TCAV und Präsentations- und Änderungsvorschläge für den Code Teil gegeben. "
$ws.Range("D15").Value = "TCAV verständnis nochmal vertieft durch die Erstellung der Präsentationsfolien, Transfer Learning, "
$ws.Range("E15").Value = "Folien für Paperpräsentation erstellt und die Gruppe unterstützt mit Vorschlägen."
$ws.Range("F15").Value = 21

# --- Row 16: (no week label / no Beitrag) ---
$ws.Range("B16").NumberFormat = "m/d/yyyy"
$ws.Range("B16").Value = 45055
$ws.Range("C16").Value = "GradCAM Paper gelesen und grob über Funktionsweise recherchiert"
$ws.Range("D16").Value = "Feature Visualisierung und Nachteile von Saliency Maps."
$ws.Range("F16").Value = 6

# --- Row 17: only the date got filled in ---
$ws.Range("B17").NumberFormat = "m/d/yyyy"
$ws.Range("B17").Value = 45062

# --- Selection / active sheet bookkeeping ---
$ws.Range("D11").Select()
$ws.Activate()
